$d = $word.ActiveDocument

# The section being removed ("Creating a database and inserting records"
# Heading3 through the last example table under "Ordering of results")
# contains four example tables plus all of their surrounding
# paragraphs/headings. Word COM's Range.Delete() silently no-ops on ranges
# that overlap a table, so drop the tables individually first.
while ($d.Tables.Count -gt 0) {
    $d.Tables.Item(1).Delete()
}

# Locate the start of the section to remove: the Heading3 paragraph
# "Създаване на база данни и вмъкване на записи" ("Creating a database and
# inserting records"), which immediately follows the last inline picture.
$findRange = $d.Content
$found = $findRange.Find.Execute("Създаване на база данни и вмъкване на записи")
$startPos = $findRange.Paragraphs.Item(1).Range.Start

# The section runs all the way to the trailing empty paragraph that sits
# right before the section break -- that paragraph must stay (only its
# formatting gets cleared below), so the deletion stops at its start.
$paras = $d.Content.Paragraphs
$endPara = $paras.Item($paras.Count)
$endPos = $endPara.Range.Start

$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()

# The remaining final (now last) paragraph used to carry only a
# "lang=bg-BG" paragraph-mark formatting; strip that so it becomes a bare
# empty paragraph, matching the target document.
$finalPara = $d.Content.Paragraphs.Last
$finalPara.Range.Select()
$word.Selection.ClearFormatting()
